$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.74000000000058"
$ws.Range("H2").Value = [double]"5.366994626498922e-10"
$ws.Range("I2").Value = [double]"5.366994626498922e-10"
$ws.Range("L2").Value = [double]"50.50856290803106"
$ws.Range("M2").Value = "[35.40122319409092, 65.6159026219712]"
$ws.Range("N2").Value = [double]"2.531835985308817e-08"
$ws.Range("O2").Value = [double]"2.531835985308817e-08"
$ws.Range("P2").Value = [double]"1.301921279866041"
$ws.Range("Q2").Value = "[0.9622896416401163, 1.641552918091965]"
$ws.Range("R2").Value = [double]"8.815312924070895e-10"
$ws.Range("S2").Value = [double]"8.815312924070895e-10"
$ws.Range("T2").Value = [double]"52.33142187738425"
$ws.Range("U2").Value = "[43.63697616769609, 61.02586758707241]"
$ws.Range("V2").Value = [double]"8.881784197001252e-16"
$ws.Range("W2").Value = [double]"8.881784197001252e-16"
$ws.Range("X2").Value = [double]"20.40648648648695"
$ws.Range("Y2").Value = [double]"19.01513513513556"
$ws.Range("Z2").Value = [double]"21.79783783783833"

# Row 3
$ws.Range("F3").Value = [double]"25.74000000000058"
$ws.Range("H3").Value = [double]"4.100215410529628e-07"
$ws.Range("I3").Value = [double]"4.100215410529628e-07"
$ws.Range("L3").Value = [double]"38.54328815858999"
$ws.Range("M3").Value = "[22.366430341809313, 54.72014597537066]"
$ws.Range("N3").Value = [double]"1.791740721079904e-05"
$ws.Range("O3").Value = [double]"1.791740721079904e-05"
$ws.Range("P3").Value = [double]"0.9622896416401163"
$ws.Range("Q3").Value = "[0.5346053564667317, 1.389973926813501]"
$ws.Range("R3").Value = [double]"4.294577757657514e-05"
$ws.Range("S3").Value = [double]"4.294577757657514e-05"
$ws.Range("T3").Value = [double]"50.81910133076473"
$ws.Range("U3").Value = "[42.34382218729834, 59.29438047423112]"
$ws.Range("V3").Value = [double]"8.881784197001252e-16"
$ws.Range("W3").Value = [double]"8.881784197001252e-16"
$ws.Range("X3").Value = [double]"21.79783783783834"
$ws.Range("Y3").Value = [double]"20.04576576576623"
$ws.Range("Z3").Value = [double]"23.54990990991045"

# Row 4
$ws.Range("F4").Value = [double]"25.74000000000058"
$ws.Range("H4").Value = [double]"3.898565914539631e-10"
$ws.Range("I4").Value = [double]"3.898565914539631e-10"
$ws.Range("L4").Value = [double]"48.25333608828575"
$ws.Range("M4").Value = "[33.41082122160226, 63.09585095496924]"
$ws.Range("N4").Value = [double]"4.784588170103632e-08"
$ws.Range("O4").Value = [double]"4.784588170103632e-08"
$ws.Range("P4").Value = [double]"0.5471843060306538"
$ws.Range("Q4").Value = "[0.2201316173686534, 0.8742369946926543]"
$ws.Range("R4").Value = [double]"0.001552041633658074"
$ws.Range("S4").Value = [double]"0.001552041633658074"
$ws.Range("T4").Value = [double]"48.75694948006571"
$ws.Range("U4").Value = "[40.94785249554242, 56.56604646458901]"
$ws.Range("V4").Value = [double]"2.220446049250313e-16"
$ws.Range("W4").Value = [double]"2.220446049250313e-16"
$ws.Range("X4").Value = [double]"23.49837837837892"
$ws.Range("Y4").Value = [double]"22.15855855855906"
$ws.Range("Z4").Value = [double]"24.83819819819877"

# Row 5
$ws.Range("F5").Value = [double]"25.74000000000058"
$ws.Range("H5").Value = [double]"2.153034861507308e-09"
$ws.Range("I5").Value = [double]"2.153034861507308e-09"
$ws.Range("L5").Value = [double]"50.18317798870523"
$ws.Range("M5").Value = "[36.31808481093984, 64.04827116647063]"
$ws.Range("N5").Value = [double]"3.794959901881612e-09"
$ws.Range("O5").Value = [double]"3.794959901881612e-09"
$ws.Range("P5").Value = [double]"0.05660527303765317"
$ws.Range("Q5").Value = "[-0.24528951649650121, 0.35850006257180755]"
$ws.Range("R5").Value = [double]"0.707470188884729"
$ws.Range("S5").Value = [double]"0.707470188884729"
$ws.Range("T5").Value = [double]"57.84067803251938"
$ws.Range("U5").Value = "[49.21244154028728, 66.46891452475148]"
$ws.Range("V5").Value = [double]"0"
$ws.Range("W5").Value = [double]"0"
$ws.Range("X5").Value = [double]"25.50810810810869"
$ws.Range("Y5").Value = [double]"24.2713513513519"
$ws.Range("Z5").Value = [double]"26.74486486486548"

# Row 6
$ws.Range("F6").Value = [double]"25.74000000000058"
$ws.Range("H6").Value = [double]"1.28864356963021e-08"
$ws.Range("I6").Value = [double]"1.28864356963021e-08"
$ws.Range("L6").Value = [double]"41.85761250404533"
$ws.Range("M6").Value = "[27.968237317806462, 55.74698769028419]"
$ws.Range("N6").Value = [double]"2.459546188937622e-07"
$ws.Range("O6").Value = [double]"2.459546188937622e-07"
$ws.Range("P6").Value = [double]"-0.3018947895341544"
$ws.Range("Q6").Value = "[-0.6666843268879239, 0.06289474781961513]"
$ws.Range("R6").Value = [double]"0.1024917190848524"
$ws.Range("S6").Value = [double]"0.1024917190848524"
$ws.Range("T6").Value = [double]"55.18228734647997"
$ws.Range("U6").Value = "[47.313726850762464, 63.050847842197484]"
$ws.Range("X6").Value = [double]"1.236756756756787"
$ws.Range("Y6").Value = [double]"-0.2576576576576606"
$ws.Range("Z6").Value = [double]"2.731171171171234"

# Row 7
$ws.Range("F7").Value = [double]"25.74000000000058"
$ws.Range("H7").Value = [double]"7.085632081071935e-11"
$ws.Range("I7").Value = [double]"7.085632081071935e-11"
$ws.Range("L7").Value = [double]"50.91236623575336"
$ws.Range("M7").Value = "[35.33042117092147, 66.49431130058525]"
$ws.Range("N7").Value = [double]"4.273584752745307e-08"
$ws.Range("O7").Value = [double]"4.273584752745307e-08"
$ws.Range("P7").Value = [double]"-0.8050527720910781"
$ws.Range("Q7").Value = "[-1.1069475616252324, -0.5031579825569237]"
$ws.Range("R7").Value = [double]"2.650989207486631e-06"
$ws.Range("S7").Value = [double]"2.650989207486631e-06"
$ws.Range("T7").Value = [double]"57.24190545364622"
$ws.Range("U7").Value = "[49.060036269055985, 65.42377463823645]"
$ws.Range("X7").Value = [double]"3.298018018018091"
$ws.Range("Y7").Value = [double]"2.061261261261306"
$ws.Range("Z7").Value = [double]"4.534774774774876"

# Row 8
$ws.Range("F8").Value = [double]"25.74000000000058"
$ws.Range("H8").Value = [double]"3.684893168376391e-09"
$ws.Range("I8").Value = [double]"3.684893168376391e-09"
$ws.Range("L8").Value = [double]"43.87115152518768"
$ws.Range("M8").Value = "[28.773134972890993, 58.969168077484376]"
$ws.Range("N8").Value = [double]"5.167211321044363e-07"
$ws.Range("O8").Value = [double]"5.167211321044363e-07"
$ws.Range("Q8").Value = "[-1.522052897234694, -0.767315923399309]"
$ws.Range("R8").Value = [double]"2.147505147043205e-07"
$ws.Range("S8").Value = [double]"2.147505147043205e-07"
$ws.Range("T8").Value = [double]"59.29321709967545"
$ws.Range("U8").Value = "[51.031285072570014, 67.55514912678088]"
$ws.Range("V8").Value = [double]"0"
$ws.Range("W8").Value = [double]"0"
$ws.Range("X8").Value = [double]"4.689369369369476"
$ws.Range("Y8").Value = [double]"3.143423423423497"
$ws.Range("Z8").Value = [double]"6.235315315315456"

# Row 9
$ws.Range("F9").Value = [double]"22"
$ws.Range("H9").Value = [double]"4.918232487938212e-11"
$ws.Range("I9").Value = [double]"4.918232487938212e-11"
$ws.Range("J9").Value = [double]"0.2987894557045725"
$ws.Range("K9").Value = [double]"0.2987894557045725"
$ws.Range("L9").Value = [double]"50.71023789241919"
$ws.Range("M9").Value = "[35.91169345815594, 65.50878232668244]"
$ws.Range("N9").Value = [double]"1.425450379954896e-08"
$ws.Range("O9").Value = [double]"1.425450379954896e-08"
$ws.Range("P9").Value = [double]"-2.138421425866927"
$ws.Range("Q9").Value = "[-2.4403162154010807, -1.8365266363327728]"
$ws.Range("R9").Value = [double]"0"
$ws.Range("S9").Value = [double]"0"
$ws.Range("T9").Value = [double]"53.70655093573375"
$ws.Range("U9").Value = "[46.046334059799044, 61.36676781166846]"
$ws.Range("X9").Value = [double]"7.487487487487488"
$ws.Range("Y9").Value = [double]"6.430430430430432"
$ws.Range("Z9").Value = [double]"8.544544544544545"

# Row 10
$ws.Range("F10").Value = [double]"22"
$ws.Range("H10").Value = [double]"1.431471940982476e-09"
$ws.Range("I10").Value = [double]"1.431471940982476e-09"
$ws.Range("J10").Value = [double]"0.8917699482542987"
$ws.Range("K10").Value = [double]"0.8917699482542987"
$ws.Range("L10").Value = [double]"45.17346580111994"
$ws.Range("M10").Value = "[32.59883591434391, 57.74809568789597]"
$ws.Range("N10").Value = [double]"4.565035904846582e-09"
$ws.Range("O10").Value = [double]"4.565035904846582e-09"
$ws.Range("P10").Value = [double]"-1.496894998106848"
$ws.Range("Q10").Value = "[-1.8113687372049245, -1.1824212590087715]"
$ws.Range("R10").Value = [double]"1.915578806688245e-12"
$ws.Range("S10").Value = [double]"1.915578806688245e-12"
$ws.Range("T10").Value = [double]"51.23545158702879"
$ws.Range("U10").Value = "[43.56133655021826, 58.909566623839325]"
$ws.Range("X10").Value = [double]"5.241241241241241"
$ws.Range("Y10").Value = [double]"4.140140140140143"
$ws.Range("Z10").Value = [double]"6.342342342342339"

# Row 11
$ws.Range("F11").Value = [double]"22"
$ws.Range("H11").Value = [double]"4.632148339345221e-07"
$ws.Range("I11").Value = [double]"4.632148339345221e-07"
$ws.Range("J11").Value = [double]"0.1167214205156615"
$ws.Range("K11").Value = [double]"0.1167214205156615"
$ws.Range("L11").Value = [double]"49.70578163377016"
$ws.Range("M11").Value = "[28.948159214824216, 70.4634040527161]"
$ws.Range("N11").Value = [double]"1.654830015906938e-05"
$ws.Range("O11").Value = [double]"1.654830015906938e-05"
$ws.Range("P11").Value = [double]"-1.081789662497386"
$ws.Range("Q11").Value = "[-1.522052897234695, -0.6415264277600778]"
$ws.Range("R11").Value = [double]"1.090183179841553e-05"
$ws.Range("S11").Value = [double]"1.090183179841553e-05"
$ws.Range("T11").Value = [double]"59.08299272464507"
$ws.Range("U11").Value = "[48.0766847925771, 70.08930065671304]"
$ws.Range("V11").Value = [double]"4.263256414560601e-14"
$ws.Range("W11").Value = [double]"4.263256414560601e-14"
$ws.Range("X11").Value = [double]"3.787787787787789"
$ws.Range("Y11").Value = [double]"2.246246246246246"
$ws.Range("Z11").Value = [double]"5.329329329329331"

Write-Output "applied all changes"